# Generate Report for Archive
# Swaps the report rows for "ae553900-8810-4298-ab6f-5195742c9be4.md" and
# "20e88495-9b5f-4072-a71e-df18e81d3aaa.md" (row 4 <-> row 5) across all
# three worksheets (Overview, zh-cn, de-de), including the hyperlink
# display text on the Overview sheet's Path column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Overview" (columns A-G)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A4").Value = "20e88495-9b5f-4072-a71e-df18e81d3aaa.md"
$ws1.Range("B4").Value = "e2e\20e88495-9b5f-4072-a71e-df18e81d3aaa.md"
$ws1.Range("C4").Value = ".md"
$ws1.Range("E4").Value = "Ready for handoff"
$ws1.Range("F4").Value = "Ready for handoff"
$ws1.Range("G4").Value = "2016-08-31 14:52:35"

$ws1.Range("A5").Value = "ae553900-8810-4298-ab6f-5195742c9be4.md"
$ws1.Range("B5").Value = "e2e\ae553900-8810-4298-ab6f-5195742c9be4.md"
$ws1.Range("C5").Value = ".md"
$ws1.Range("E5").Value = "In Translation"
$ws1.Range("F5").Value = "In Translation"
$ws1.Range("G5").Value = "2016-08-31 14:51:47"

# Rebuild hyperlinks on column B (Path And Name) in the new row order.
# Engine API only exposes whole-collection Delete(); re-add every link.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c447ec244f689c90f080331e913f5886c7b69374/e2e/70bbdfec-cd28-4f44-920e-fb9ddac3537f.md", "", "", "e2e\70bbdfec-cd28-4f44-920e-fb9ddac3537f.md")
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8794404c719ca6cb1bbc14ff808f59b806600140/e2e/103e825a-ebd3-433f-8009-aaedf32cae49.md", "", "", "e2e\103e825a-ebd3-433f-8009-aaedf32cae49.md")
$ws1.Hyperlinks.Add($ws1.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8794404c719ca6cb1bbc14ff808f59b806600140/e2e/ae553900-8810-4298-ab6f-5195742c9be4.md", "", "", "e2e\20e88495-9b5f-4072-a71e-df18e81d3aaa.md")
$ws1.Hyperlinks.Add($ws1.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5e2a09039fe4dc271e7d0205b7c4411f5148fb6c/e2e/20e88495-9b5f-4072-a71e-df18e81d3aaa.md", "", "", "e2e\ae553900-8810-4298-ab6f-5195742c9be4.md")
$ws1.Hyperlinks.Add($ws1.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/04b0623ccd1c951b7d867c8917073ac2d382c4af/e2e/a02783f8-0796-4b0c-bed2-1b4f1c5eb63a.md", "", "", "e2e\a02783f8-0796-4b0c-bed2-1b4f1c5eb63a.md")
$ws1.Hyperlinks.Add($ws1.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/64188b06fcf46896a066d032512cba7e1dffbbba/e2e/b2121808-9ac1-4fd6-a3d6-22fe1d966b9c.md", "", "", "e2e\b2121808-9ac1-4fd6-a3d6-22fe1d966b9c.md")

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn" (columns A-P)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A4").Value = "20e88495-9b5f-4072-a71e-df18e81d3aaa.md"
$ws2.Range("C4").Value = "Ready for handoff"
$ws2.Range("G4").Value = "20e88495-9b5f-4072-a71e-df18e81d3aaa.57ac5ae7e920a304d65acd4e0dc0f1b7de521c12.zh-cn.xlf"
$ws2.Range("H4").Value = "2016-08-31 14:52:29"

$ws2.Range("A5").Value = "ae553900-8810-4298-ab6f-5195742c9be4.md"
$ws2.Range("C5").Value = "In Translation"
$ws2.Range("G5").Value = "ae553900-8810-4298-ab6f-5195742c9be4.ebacd5ec84660e7d6f03385cb5bb23f0cb71b73d.zh-cn.xlf"
$ws2.Range("H5").Value = "2016-08-31 14:51:43"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c447ec244f689c90f080331e913f5886c7b69374/e2e/70bbdfec-cd28-4f44-920e-fb9ddac3537f.md", "", "", "70bbdfec-cd28-4f44-920e-fb9ddac3537f.md")
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/bb9f9c0551e1e277220c173d62a0da33966f0d9d/e2e/70bbdfec-cd28-4f44-920e-fb9ddac3537f.md", "", "", "70bbdfec-cd28-4f44-920e-fb9ddac3537f.md")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8794404c719ca6cb1bbc14ff808f59b806600140/e2e/103e825a-ebd3-433f-8009-aaedf32cae49.md", "", "", "103e825a-ebd3-433f-8009-aaedf32cae49.md")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8794404c719ca6cb1bbc14ff808f59b806600140/e2e/ae553900-8810-4298-ab6f-5195742c9be4.md", "", "", "20e88495-9b5f-4072-a71e-df18e81d3aaa.md")
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5e2a09039fe4dc271e7d0205b7c4411f5148fb6c/e2e/20e88495-9b5f-4072-a71e-df18e81d3aaa.md", "", "", "ae553900-8810-4298-ab6f-5195742c9be4.md")
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/04b0623ccd1c951b7d867c8917073ac2d382c4af/e2e/a02783f8-0796-4b0c-bed2-1b4f1c5eb63a.md", "", "", "a02783f8-0796-4b0c-bed2-1b4f1c5eb63a.md")
$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/64188b06fcf46896a066d032512cba7e1dffbbba/e2e/b2121808-9ac1-4fd6-a3d6-22fe1d966b9c.md", "", "", "b2121808-9ac1-4fd6-a3d6-22fe1d966b9c.md")

# ---------------------------------------------------------------------
# Sheet 3: "de-de" (columns A-P)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A4").Value = "20e88495-9b5f-4072-a71e-df18e81d3aaa.md"
$ws3.Range("C4").Value = "Ready for handoff"
$ws3.Range("G4").Value = "20e88495-9b5f-4072-a71e-df18e81d3aaa.57ac5ae7e920a304d65acd4e0dc0f1b7de521c12.de-de.xlf"
$ws3.Range("H4").Value = "2016-08-31 14:52:35"

$ws3.Range("A5").Value = "ae553900-8810-4298-ab6f-5195742c9be4.md"
$ws3.Range("C5").Value = "In Translation"
$ws3.Range("G5").Value = "ae553900-8810-4298-ab6f-5195742c9be4.ebacd5ec84660e7d6f03385cb5bb23f0cb71b73d.de-de.xlf"
$ws3.Range("H5").Value = "2016-08-31 14:51:47"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c447ec244f689c90f080331e913f5886c7b69374/e2e/70bbdfec-cd28-4f44-920e-fb9ddac3537f.md", "", "", "70bbdfec-cd28-4f44-920e-fb9ddac3537f.md")
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4fc68cbaafcaf47c01aa618d712e6e3b07cb7e90/e2e/70bbdfec-cd28-4f44-920e-fb9ddac3537f.md", "", "", "70bbdfec-cd28-4f44-920e-fb9ddac3537f.md")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8794404c719ca6cb1bbc14ff808f59b806600140/e2e/103e825a-ebd3-433f-8009-aaedf32cae49.md", "", "", "103e825a-ebd3-433f-8009-aaedf32cae49.md")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8794404c719ca6cb1bbc14ff808f59b806600140/e2e/ae553900-8810-4298-ab6f-5195742c9be4.md", "", "", "20e88495-9b5f-4072-a71e-df18e81d3aaa.md")
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5e2a09039fe4dc271e7d0205b7c4411f5148fb6c/e2e/20e88495-9b5f-4072-a71e-df18e81d3aaa.md", "", "", "ae553900-8810-4298-ab6f-5195742c9be4.md")
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/04b0623ccd1c951b7d867c8917073ac2d382c4af/e2e/a02783f8-0796-4b0c-bed2-1b4f1c5eb63a.md", "", "", "a02783f8-0796-4b0c-bed2-1b4f1c5eb63a.md")
$ws3.Hyperlinks.Add($ws3.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/64188b06fcf46896a066d032512cba7e1dffbbba/e2e/b2121808-9ac1-4fd6-a3d6-22fe1d966b9c.md", "", "", "b2121808-9ac1-4fd6-a3d6-22fe1d966b9c.md")
